# Weekly data refresh: insert the newest week's row of Mandarina price data
# at row 25 (just below the static header block), pushing all the existing
# historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 25 - this shifts rows 25:107 down to
# 26:108 (and carries their formatting/styles with them, Excel-native style).
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with this week's entry.
$ws.Range("A25").Value = 11
$ws.Range("B25").Value = "Vega Monumental Concepción"
$ws.Range("C25").Value = "Bíobío"
$ws.Range("D25").Value = 44701
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100102
$ws.Range("H25").Value = "Cítricos"
$ws.Range("I25").Value = 100102004
$ws.Range("J25").Value = "Mandarina"
$ws.Range("K25").Value = "Clemenuless"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 15000
$ws.Range("P25").Value = 14500
$ws.Range("Q25").Value = "$/bandeja 10 kilos"
$ws.Range("R25").Value = "Provincia de Limarí"
$ws.Range("S25").Value = 1450
$ws.Range("T25").Value = 10
